$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NET")

# Update row 7 ("Change in payables and accrued liability") values B7:F7
$ws.Range("B7").Value = 62000000.0
$ws.Range("C7").Value = 59000000.0
$ws.Range("D7").Value = 46732000.0
$ws.Range("E7").Value = 29630000.0
$ws.Range("F7").Value = 16652000.0
